$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.675.24'
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').Value = '2.031.36'
$ws.Range('E3').Value = '  -0.12%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '227.43'
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('E6').Value = '  -1.07%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '59.73'
$ws.Range('E7').Value = '  -0.98%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('E9').Value = '  -1.68%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0826'
$ws.Range('E10').Value = '  +1.79%  '
$ws.Range('D12').Value = '2.332.64'
$ws.Range('E12').Value = '  -0.21%  '
$ws.Range('E13').Value = '  -0.93%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '20.99'
$ws.Range('E14').Value = '  -0.22%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.773'
$ws.Range('E15').Value = '  +2.32%  '
$ws.Range('E16').Value = '  +0.60%  '
$ws.Range('D17').Value = '2.028.79'
$ws.Range('E17').Value = '  -0.74%  '
$ws.Range('D18').Value = '37.646.82'
$ws.Range('E18').Value = '  -0.49%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '5.93'
$ws.Range('E19').Value = '  -2.74%  '
$ws.Range('B20').Value = 'Litecoin'
$ws.Range('C20').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '69.34'
$ws.Range('E20').Value = '  -0.52%  '
$ws.Range('D21').Value = '0.0₃0819'
$ws.Range('E21').Value = '  -0.81%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '223.70'
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('E24').Value = '  -0.28%  '
$ws.Range('E25').Value = '  +3.98%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '167.62'
$ws.Range('E26').Value = '  +0.88%  '
$ws.Range('E27').Value = '  +1.79%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.127'
$ws.Range('E28').Value = '  -1.62%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '18.70'
$ws.Range('E29').Value = '  -0.91%  '
$ws.Range('E30').Value = '  -2.26%  '
$ws.Range('E31').Value = '  +0.31%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.22'
$ws.Range('E32').Value = '  +8.26%  '
$ws.Range('E33').Value = '  -1.83%  '
$ws.Range('E34').Value = '  +0.12%  '
$ws.Range('E35').Value = '  -0.96%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.50'
$ws.Range('E36').Value = '  +3.41%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.34'
$ws.Range('E37').Value = '  +3.01%  '
$ws.Range('E38').Value = '  +5.13%  '
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '17.84'
$ws.Range('E40').Value = '  +8.04%  '
$ws.Range('D41').Value = '1.520.12'
$ws.Range('E41').Value = '  -1.09%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '97.02'
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('E43').Value = '  -1.43%  '
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0906'
$ws.Range('E45').Value = '  -1.68%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.14'
$ws.Range('E46').Value = '  +2.91%  '
$ws.Range('E47').Value = '  -0.94%  '
$ws.Range('E48').Value = '  -0.04%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '7.08'
$ws.Range('E49').Value = '  +0.47%  '
$ws.Range('B50').Value = 'MXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.93'
$ws.Range('E50').Value = '  -1.50%  '
$ws.Range('D51').Value = '2.223.05'
$ws.Range('E51').Value = '  -0.21%  '
